$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers -----------------------------------------------------
# Preserve the existing "Average" shared string by writing it to F1 first,
# then overwrite E1 with the new "Day 4" header (so shared-string indices
# come out in the same order as the target file).
$ws.Range("F1").Value = "Average"
$ws.Range("E1").Value = "Day 4 (2023-02-25)"

# --- Row 2 ---------------------------------------------------------------
# Establish the "0.0" number format first (cellXfs index 1 in the target).
$ws.Range("E2").Value = 20.8
$ws.Range("F2").Formula = "=AVERAGE(B2:E2)"
$ws.Range("F2").NumberFormat = "0.0"
$ws.Range("G2").Formula = "=F2*0.8"
$ws.Range("G2").NumberFormat = "0.0"

# Ratio cell used by the new "G" column (displayed as a percentage,
# cellXfs index 2 in the target).
$ws.Range("G1").Value = 0.8
$ws.Range("G1").NumberFormat = "0%"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("E3").Value = 20.4
$ws.Range("F3").Formula = "=AVERAGE(B3:E3)"
$ws.Range("F3").NumberFormat = "0.0"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("E4").Value = 23.5
$ws.Range("F4").Formula = "=AVERAGE(B4:E4)"
$ws.Range("F4").NumberFormat = "0.0"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("E5").Value = 25.1
$ws.Range("F5").Formula = "=AVERAGE(B5:E5)"
$ws.Range("F5").NumberFormat = "0.0"

# G3:G5 were filled together as one shared formula (mirrors the original
# E3:E5 shared-formula group), so create them the same way.
$ws.Range("G3:G5").Formula = "=F3*0.8"
$ws.Range("G3").NumberFormat = "0.0"
$ws.Range("G5").NumberFormat = "0.0"
# Evaluating the shared formula on G4 picks up F4's "0.0" number format as a
# side effect (same-row precedent cell) - reset G4 back to the default style
# to match the target (G4 keeps General formatting).
$ws.Range("G4").Style = "Normal"

# --- Row 6: blank separator row with the next week's headers ------------
$ws.Range("A7").Value = "coop014"
$ws.Range("A8").Value = "coop015"
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = "Day 1 (2023-03-02)"
$ws.Range("C6").Value = "Day 2 (2023-03-03)"
$ws.Range("D6").Value = "Day 3 (2023-03-04)"
$ws.Range("F6").NumberFormat = "0.0"
$ws.Range("G6").NumberFormat = "0.0"

# --- Row 7: new pair (coop014) -------------------------------------------
$ws.Range("B7").Value = 27.8
$ws.Range("C7").Value = 28.1

# --- Row 8: new pair (coop015) -------------------------------------------
$ws.Range("B8").Value = 29.2
$ws.Range("C8").Value = 29.1

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 16.1666667
$ws.Columns.Item(5).ColumnWidth = 16.1666667

# --- Selection -------------------------------------------------------------
$ws.Range("D7").Select()
